$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1013
$ws1.Range("F9").Value = 1480
$ws1.Range("F13").Value = 499
$ws1.Range("F14").Value = 1677
$ws1.Range("F16").Value = 817
$ws1.Range("F21").Value = 1151
$ws1.Range("F23").Value = 411
$ws1.Range("F24").Value = 27
$ws1.Range("F25").Value = 3571
$ws1.Range("F28").Value = 1577

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 36

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 36
$ws4.Range("F16").Value = 1013
$ws4.Range("F20").Value = 1480
$ws4.Range("F24").Value = 499
$ws4.Range("F25").Value = 1677
$ws4.Range("F27").Value = 817
$ws4.Range("F34").Value = 1151
$ws4.Range("F36").Value = 411
$ws4.Range("F37").Value = 27
$ws4.Range("F38").Value = 3571
$ws4.Range("F41").Value = 1577
